$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\aa\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\.DS_Store', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\aa\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\xdd\Docs\.DS_Store', $false, '.DS_Store', '.DS_Store'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\aa\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\.DS_Store', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\aa\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\xdd\Docs\.DS_Store', $false, '.DS_Store', '.DS_Store'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\aa\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\.DS_Store', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\aa\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\xdd\bin\.DS_Store', $false, '.DS_Store', '.DS_Store'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\aa\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\.DS_Store', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\aa\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\xdd\bin\.DS_Store', $false, '.DS_Store', '.DS_Store'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\08456059.pdf', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\08456059.pdf', $true, '08456059.pdf', '08456059.pdf'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\20180006549.pdf', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\20180006549.pdf', $true, '20180006549.pdf', '20180006549.pdf'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\235434231_RAW_PAN.pix', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\235434231_RAW_PAN.pix', $true, '235434231_RAW_PAN.pix', '235434231_RAW_PAN.pix'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\235434251_RAW_PAN_LOG.txt', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\235434251_RAW_PAN_LOG.txt', $true, '235434251_RAW_PAN_LOG.txt', '235434251_RAW_PAN_LOG.txt'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\3.205-CVE_to_be_installed.txt', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\3.205-CVE_to_be_installed.txt', $true, '3.205-CVE_to_be_installed.txt', '3.205-CVE_to_be_installed.txt'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\3.207-CVE_to_be_installed.txt', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\3.207-CVE_to_be_installed.txt', $true, '3.207-CVE_to_be_installed.txt', '3.207-CVE_to_be_installed.txt'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\anaconda-ks.cfg', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\anaconda-ks.cfg', $true, 'anaconda-ks.cfg', 'anaconda-ks.cfg'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\authorized_keys-4-8-17', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\authorized_keys-4-8-17', $true, 'authorized_keys-4-8-17', 'authorized_keys-4-8-17'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\BAND_RPC.txt', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\BAND_RPC.txt', $true, 'BAND_RPC.txt', 'BAND_RPC.txt'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\buttons.print.min.js', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\buttons.print.min.js', $true, 'buttons.print.min.js', 'buttons.print.min.js'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\catalina.2023-07-14.log', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\catalina.2023-07-14.log', $true, 'catalina.2023-07-14.log', 'catalina.2023-07-14.log'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\catalina.2023-07-16.log', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\catalina.2023-07-16.log', $true, 'catalina.2023-07-16.log', 'catalina.2023-07-16.log'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\catalina.2023-07-23.log', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\catalina.2023-07-23.log', $true, 'catalina.2023-07-23.log', 'catalina.2023-07-23.log'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\catalina.2023-08-26.log', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\catalina.2023-08-26.log', $true, 'catalina.2023-08-26.log', 'catalina.2023-08-26.log'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\covid_excel_20-38.xlsx', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\covid_excel_20-38.xlsx', $true, 'covid_excel_20-38.xlsx', 'covid_excel_20-38.xlsx'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\ddn-restapi-core-gs-4.2.1-3.3d53e4_dirty.noarch.rpm', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\ddn-restapi-core-gs-4.2.1-3.3d53e4_dirty.noarch.rpm', $true, 'ddn-restapi-core-gs-4.2.1-3.3d53e4_dirty.noarch.rpm', 'ddn-restapi-core-gs-4.2.1-3.3d53e4_dirty.noarch.rpm'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\EaseUS Partition Master 9.2.2 Technician Edition (FULL + Patch).zip', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\EaseUS Partition Master 9.2.2 Technician Edition (FULL + Patch).zip', $true, 'EaseUS Partition Master 9.2.2 Technician Edition (FULL + Patch).zip', 'EaseUS Partition Master 9.2.2 Technician Edition (FULL + Patch).zip'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\fio-3.1-1.el6.x86_64.rpm', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\fio-3.1-1.el6.x86_64.rpm', $true, 'fio-3.1-1.el6.x86_64.rpm', 'fio-3.1-1.el6.x86_64.rpm'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\gpfs.gss.pmsensors-4.2.3-5.el7.x86_64.rpm', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\gpfs.gss.pmsensors-4.2.3-5.el7.x86_64.rpm', $true, 'gpfs.gss.pmsensors-4.2.3-5.el7.x86_64.rpm', 'gpfs.gss.pmsensors-4.2.3-5.el7.x86_64.rpm'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\gwm66c08.jpg', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\gwm66c08.jpg', $true, 'gwm66c08.jpg', 'gwm66c08.jpg'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\import CSV.txt', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\import CSV.txt', $true, 'import CSV.txt', 'import CSV.txt'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\KKD geo tagging (1).xlsx', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\KKD geo tagging (1).xlsx', $true, 'KKD geo tagging (1).xlsx', 'KKD geo tagging (1).xlsx'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\LS____.pdf', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\LS____.pdf', $true, 'LS____.pdf', 'LS____.pdf'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\mmhealth.sh', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\mmhealth.sh', $true, 'mmhealth.sh', 'mmhealth.sh'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\mounted', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\mounted', $true, 'mounted', 'mounted'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\NewGCPs.bk', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\NewGCPs.bk', $true, 'NewGCPs.bk', 'NewGCPs.bk'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\SR159360_ddn_showall_20210222-13_54_51.tar.gz', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\SR159360_ddn_showall_20210222-13_54_51.tar.gz', $true, 'SR159360_ddn_showall_20210222-13_54_51.tar.gz', 'SR159360_ddn_showall_20210222-13_54_51.tar.gz'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\SR170158_ddn_showall_20210902-13_20_38.tar.gz', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\SR170158_ddn_showall_20210902-13_20_38.tar.gz', $true, 'SR170158_ddn_showall_20210902-13_20_38.tar.gz', 'SR170158_ddn_showall_20210902-13_20_38.tar.gz'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\syslog', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\syslog', $true, 'syslog', 'syslog'),
    @('D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-1\table.sql', 'D:\\adhvik\\adh\\Hackathon\\space hack\\top12 data\\topic12\\dataset1\folder-2\table.sql', $true, 'table.sql', 'table.sql')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$ws.Range("A36:E47").ClearContents() | Out-Null

Write-Output $ws.UsedRange.Address()
